$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N, shifting "Late" / "Outstanding"
# (and their data) one column to the right.
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab and select S5, matching the
# new active tab + selection recorded in the workbook.
$ws.Activate()
$ws.Range("S5").Select()
